$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row after the last data row (row 60), pushing the blank
#    gap rows and the signature block down by one.
$ws.Rows("61:61").Insert()

# 2. Populate the new row 61 with a new period record (2509), copying the
#    repeated worker info from the row above it (row 60) and using the same
#    totals as that row.
$ws.Range("B61").Value = $ws.Range("B60").Value2
$ws.Range("C61").Value = $ws.Range("C60").Value2
$ws.Range("D61").Value = $ws.Range("D60").Value2
$ws.Range("E61").Value = "2509"
$ws.Range("F61").Value = 164658
$ws.Range("G61").Value = 4116447

# 3. Give the new last row the "closing" bottom-border styling that row 60
#    used to have, and restore row 60 to the regular interior-row styling
#    (matching row 59) now that it is no longer the last row.
$ws.Range("B60:J60").Borders.Item(9).LineStyle = 1
$ws.Range("B60:J60").Borders.Item(9).Weight = 2
$ws.Range("B61:J61").Borders.Item(9).LineStyle = 1
$ws.Range("B61:J61").Borders.Item(9).Weight = 2

# 4. Update the totals that changed because of the new period record.
$ws.Range("E11").Value = 6379399
$ws.Range("F13").Value = 44
